$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.572976999999999
$ws.Range("H2").Value = 22.718931
$ws.Range("I2").Value = 0.4497670593913077
$ws.Range("J2").Value = 0.4497670593913078
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 146.0459156666667
$ws.Range("N2").Value = 438.137747
$ws.Range("O2").Value = 0.4862506770104965
$ws.Range("P2").Value = 0.4862506770104965
$ws.Range("Q2").Value = 1106.002360287606
$ws.Range("R2").Value = 9954.021242588457
$ws.Range("S2").Value = 0.2186995371260436
$ws.Range("T2").Value = 0.2186995371260436
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.572976999999999
$ws.Range("H3").Value = 22.718931
$ws.Range("I3").Value = 0.4497670593913077
$ws.Range("J3").Value = 0.4497670593913078
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 122.7232436666666
$ws.Range("N3").Value = 368.169731
$ws.Range("O3").Value = 0.408599309644787
$ws.Range("P3").Value = 0.408599309644787
$ws.Range("Q3").Value = 929.380301653062
$ws.Range("R3").Value = 8364.422714877559
$ws.Range("S3").Value = 0.1837745099682542
$ws.Range("T3").Value = 0.1837745099682543
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.572976999999999
$ws.Range("H4").Value = 22.718931
$ws.Range("I4").Value = 0.4497670593913077
$ws.Range("J4").Value = 0.4497670593913078
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 31.58192
$ws.Range("N4").Value = 94.74576
$ws.Range("O4").Value = 0.1051500133447165
$ws.Range("P4").Value = 0.1051500133447165
$ws.Range("Q4").Value = 239.16915377584
$ws.Range("R4").Value = 2152.52238398256
$ws.Range("S4").Value = 0.0472930122970099
$ws.Range("T4").Value = 0.04729301229700991
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.784025666666667
$ws.Range("H5").Value = 8.352077
$ws.Range("I5").Value = 0.1653462089435359
$ws.Range("J5").Value = 0.1653462089435359
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 146.0459156666667
$ws.Range("N5").Value = 438.137747
$ws.Range("O5").Value = 0.4862506770104965
$ws.Range("P5").Value = 0.4862506770104965
$ws.Range("Q5").Value = 406.5955777278355
$ws.Range("R5").Value = 3659.360199550519
$ws.Range("S5").Value = 0.08039970603991335
$ws.Range("T5").Value = 0.08039970603991335
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.784025666666667
$ws.Range("H6").Value = 8.352077
$ws.Range("I6").Value = 0.1653462089435359
$ws.Range("J6").Value = 0.1653462089435359
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 122.7232436666666
$ws.Range("N6").Value = 368.169731
$ws.Range("O6").Value = 0.408599309644787
$ws.Range("P6").Value = 0.408599309644787
$ws.Range("Q6").Value = 341.6646602645874
$ws.Range("R6").Value = 3074.981942381286
$ws.Range("S6").Value = 0.06756034682671148
$ws.Range("T6").Value = 0.0675603468267115
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.784025666666667
$ws.Range("H7").Value = 8.352077
$ws.Range("I7").Value = 0.1653462089435359
$ws.Range("J7").Value = 0.1653462089435359
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 31.58192
$ws.Range("N7").Value = 94.74576
$ws.Range("O7").Value = 0.1051500133447165
$ws.Range("P7").Value = 0.1051500133447165
$ws.Range("Q7").Value = 87.92487588261334
$ws.Range("R7").Value = 791.3238829435199
$ws.Range("S7").Value = 0.01738615607691108
$ws.Range("T7").Value = 0.01738615607691108
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.480551000000001
$ws.Range("H8").Value = 19.441653
$ws.Range("I8").Value = 0.3848867316651562
$ws.Range("J8").Value = 0.3848867316651562
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 146.0459156666667
$ws.Range("N8").Value = 438.137747
$ws.Range("O8").Value = 0.4862506770104965
$ws.Range("P8").Value = 0.4862506770104965
$ws.Range("Q8").Value = 946.4580048195326
$ws.Range("R8").Value = 8518.122043375792
$ws.Range("S8").Value = 0.1871514338445395
$ws.Range("T8").Value = 0.1871514338445395
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.480551000000001
$ws.Range("H9").Value = 19.441653
$ws.Range("I9").Value = 0.3848867316651562
$ws.Range("J9").Value = 0.3848867316651562
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 122.7232436666666
$ws.Range("N9").Value = 368.169731
$ws.Range("O9").Value = 0.408599309644787
$ws.Range("P9").Value = 0.408599309644787
$ws.Range("Q9").Value = 795.3142394672603
$ws.Range("R9").Value = 7157.828155205343
$ws.Range("S9").Value = 0.1572644528498212
$ws.Range("T9").Value = 0.1572644528498212
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.480551000000001
$ws.Range("H10").Value = 19.441653
$ws.Range("I10").Value = 0.3848867316651562
$ws.Range("J10").Value = 0.3848867316651562
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.58192
$ws.Range("N10").Value = 94.74576
$ws.Range("O10").Value = 0.1051500133447165
$ws.Range("P10").Value = 0.1051500133447165
$ws.Range("Q10").Value = 204.66824323792
$ws.Range("R10").Value = 1842.01418914128
$ws.Range("S10").Value = 0.0472930122970099
$ws.Range("T10").Value = 0.04729301229700991
